$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spelling of "Mattew Hudson-Smith" -> "Matthew Hudson-Smith" in C6
$ws.Range("C6").Value = "Matthew Hudson-Smith"

# Update the selected cell/active cell to C7
$ws.Range("C7").Select()
